$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Framework 7.1 update: new Teradata datasource ("TEWSA_APP" instead of
# "TEWSA_OWN") with a new host/db name; username & password cells both
# point at the new account name, the JDBC url moves to the new server.
$ws.Range("E2").Value = "jdbc:teradata://STDCAS02-1-2.sede.corp.sanpaoloimi.com/DATABASE=TEWSA0W"
$ws.Range("F2").Value = "TEWSA_APP"
$ws.Range("G2").Value = "TEWSA_APP"

# Duplicate the trailing decorative AutoShape (a borderless white
# rectangle drawn behind the cell comments) so the sheet keeps one
# extra copy of it, same size/position as the existing one
# (5429250 x 9515475 EMU == 427.5 x 749.25 points).
$lastShape = $ws.Shapes.Item($ws.Shapes.Count)
$newShape = $ws.Shapes.AddShape(1, $lastShape.Left, $lastShape.Top, 427.5, 749.25)
$newShape.Name = $lastShape.Name
$newShape.Fill.ForeColor.RGB = 0xFFFFFF
$newShape.Line.ForeColor.RGB = 0x000000

$ws.Range("G2").Select()
